# Update the "取得日時" (acquisition timestamp) column on the "ランサーズ" sheet
# for all data rows (2 through 20) to the new timestamp value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-15 18:37:47"

for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
